$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.289.87'
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").Value = '2.528.46'
$ws.Range("E3").Value = '  +1.76%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '536.34'
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.72'
$ws.Range("E6").Value = '  -2.61%  '
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.564'
$ws.Range("E8").Value = '  -1.65%  '
$ws.Range("D9").Value = '2.536.10'
$ws.Range("E9").Value = '  +0.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0995'
$ws.Range("E10").Value = '  -0.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.161'
$ws.Range("E11").Value = '  +1.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.40'
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.355'
$ws.Range("E13").Value = '  +0.79%  '
$ws.Range("D14").Value = '2.976.79'
$ws.Range("E14").Value = '  +1.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.19'
$ws.Range("E15").Value = '  -2.70%  '
$ws.Range("D16").Value = '59.280.91'
$ws.Range("E16").Value = '  +0.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000141'
$ws.Range("E17").Value = '  +1.44%  '
$ws.Range("D18").Value = '2.535.84'
$ws.Range("E18").Value = '  +1.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.99'
$ws.Range("E19").Value = '  -2.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.24'
$ws.Range("E20").Value = '  -0.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.40'
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.82'
$ws.Range("E23").Value = '  +1.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.04'
$ws.Range("E24").Value = '  +2.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.419'
$ws.Range("E25").Value = '  -4.31%  '
$ws.Range("E26").Value = '  +3.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.83'
$ws.Range("E28").Value = '  +1.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.77'
$ws.Range("E29").Value = '  -0.35%  '
$ws.Range("D30").Value = '0.0₃0773'
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("E31").Value = '  +0.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '161.54'
$ws.Range("E32").Value = '  +1.25%  '
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.14'
$ws.Range("E34").Value = '  -8.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.45'
$ws.Range("E35").Value = '  +0.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.50'
$ws.Range("E36").Value = '  -0.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.23'
$ws.Range("E37").Value = '  -4.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.59'
$ws.Range("E38").Value = '  -1.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.98'
$ws.Range("E39").Value = '  +0.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.66'
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.34'
$ws.Range("E41").Value = '  -8.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '287.96'
$ws.Range("E42").Value = '  -6.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.804'
$ws.Range("E43").Value = '  -1.61%  '
$ws.Range("E44").Value = '  +0.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.602'
$ws.Range("E45").Value = '  +1.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.85'
$ws.Range("E46").Value = '  +0.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '124.53'
$ws.Range("E47").Value = '  +0.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0929'
$ws.Range("E48").Value = '  +0.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.62'
$ws.Range("E49").Value = '  +0.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0510'
$ws.Range("E50").Value = '  -0.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0223'
$ws.Range("E51").Value = '  -1.92%  '
